# Apply the "Data" -> "Summary" Paraguay MSME summary update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the worksheet from "Data" to "Summary"
# ---------------------------------------------------------------------------
$ws.Name = "Summary"

# Re-assert the formatting of the two untouched header cells: the COM
# round-trip of this workbook's original named cell styles does not always
# keep the font linked to the right cellXfs entry, so explicitly re-apply
# the intended direct formatting to guarantee it survives the save.
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2. Clear the cells whose old content is not reused at the same address in
#    the new layout (their data is being relocated further down the sheet).
# ---------------------------------------------------------------------------
$cellsToClear = @("B5","C5","D5","A6","D6","D7","A8","D8","A9")
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Clear()
}

# ---------------------------------------------------------------------------
# 3. Write the new / relocated content into its final position.
# ---------------------------------------------------------------------------

# New bold+underlined "title_" style line describing the source type.
$ws.Range("A7").Value2 = "Source Type: Statistical Institution (Most Widely Used)"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

# Column headers (moved from row 5 to row 9)
$ws.Range("B9").Value2 = "Micro"
$ws.Range("B9").Font.Bold = $true
$ws.Range("C9").Value2 = "SMEs"
$ws.Range("C9").Font.Bold = $true
$ws.Range("D9").Value2 = "MSMEs"
$ws.Range("D9").Font.Bold = $true

# New "Employment (% of total)" row
$ws.Range("A10").Value2 = "Employment (% of total)"
$ws.Range("A10").Font.Bold = $true
$ws.Range("D10").Value2 = "'61.6"
$ws.Range("D10").Style = "Normal"

# "Enterprises (absolute #)" row (moved from row 6 to row 11)
$ws.Range("A11").Value2 = "Enterprises (absolute #)"
$ws.Range("A11").Font.Bold = $true
$ws.Range("D11").Value2 = "'217250"
$ws.Range("D11").Style = "Normal"

# "Enterprises density (per 1000 people)" row (moved from row 7 to row 12)
$ws.Range("A12").Value2 = "Enterprises density (per 1000 people)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("D12").Value2 = "'33.6"
$ws.Range("D12").Style = "Normal"

# "Employment (absolute #)" row (moved from row 8 to row 13)
$ws.Range("A13").Value2 = "Employment (absolute #)"
$ws.Range("A13").Font.Bold = $true
$ws.Range("D13").Value2 = "'492181"
$ws.Range("D13").Style = "Normal"

# "Enterprises (% of total)" row (moved from row 9 to row 14)
$ws.Range("A14").Value2 = "Enterprises (% of total)"
$ws.Range("A14").Font.Bold = $true
$ws.Range("D14").Value2 = "'96.9"
$ws.Range("D14").Style = "Normal"

# "Source: DGEEC, 2010" (moved from row 10 to row 15), italic "source" style
$ws.Range("A15").Value2 = "Source: DGEEC, 2010"
$ws.Range("A15").Font.Italic = $true

# New "DGEEC" bold title line
$ws.Range("A23").Value2 = "DGEEC"
$ws.Range("A23").Font.Bold = $true

# New italic citation / source detail line
$ws.Range("A24").Value2 = 'Dirección General de Estadística, Encuestas y Censos (DGEEC), "Censo Económico Nacional 2011", 2013, p.57. Available at http://www.dgeec.gov.py/Publicaciones/Biblioteca/CEN2011/resultados_finales_CEN.pdf'
$ws.Range("A24").Font.Italic = $true
